$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Bring row 26/27's hour/description cells into line with the rest of
#     the table (rows 14-25), which no longer carry the (visually
#     redundant) "applyFill" flavoured styles. Copy the format from row 25
#     (B25/C25) down onto B26:C27.
$ws.Range("B25:C25").Copy()
$ws.Range("B26:C27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Add the new work-log entry (row 28) ---
$newRow = 28
$ws.Range("A$newRow").Value = 45718
$ws.Range("B$newRow").Value = 4
$ws.Range("C$newRow").Value = "Explored more tools and commands (little bit more on commix, fierce, Dirb, WPScan, ettercap, xsser)"

# Match the formatting already used throughout the table (same as row 25/26/27)
$ws.Range("A25:C25").Copy()
$ws.Range("A$newRow`:C$newRow").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Mirror the cursor position left behind after typing the new row ---
$ws.Range("D$newRow").Select()
